$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '71.007.12'
$ws.Range("E2").Value = '  +0.58%  '

# Row 3
$ws.Range("D3").Value = '3.543.94'
$ws.Range("E3").Value = '  -0.60%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.56%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.27'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.28%  '

# Row 7
$ws.Range("D7").Value = '3.537.91'
$ws.Range("E7").Value = '  -0.71%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.613'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.05%  '

# Row 9
$ws.Range("E9").Value = '  +0.06%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.199'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.58%  '

# Row 11
$ws.Range("E11").Value = '  -3.58%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.589'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.42%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.72'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.04%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000277'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.13%  '

# Row 15
$ws.Range("D15").Value = '4.115.94'
$ws.Range("E15").Value = '  -0.57%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.72%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '614.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.50%  '

# Row 18
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.545.14'
$ws.Range("E18").Value = '  -0.60%  '

# Row 19
$ws.Range("B19").Value = 'WrappedBTC'
$ws.Range("C19").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D19").Value = '71.102.78'
$ws.Range("E19").Value = '  +0.55%  '

# Row 20
$ws.Range("E20").Value = '  +1.36%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.78'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.30%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.890'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.69%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.53%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.33%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.40%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.79'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.11%  '

# Row 27
$ws.Range("E27").Value = '  -0.01%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.23%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '33.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.04%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.00%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.06'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.11%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.18'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.93%  '

# Row 33
$ws.Range("E33").Value = '  -0.25%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.87'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.84%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '634.23'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +10.03%  '

# Row 36
$ws.Range("E36").Value = '  -0.72%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.87'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.29%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.51'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.71%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0479'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.24%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.04'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.56%  '

# Row 41
$ws.Range("E41").Value = '  +2.58%  '

# Row 42
$ws.Range("E42").Value = '  +0.03%  '

# Row 43
$ws.Range("D43").Value = '0.0₃0745'
$ws.Range("E43").Value = '  +5.55%  '

# Row 44
$ws.Range("D44").Value = '3.373.96'
$ws.Range("E44").Value = '  -0.62%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.98'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.44%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.315'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.99%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.18%  '

# Row 48
$ws.Range("E48").Value = '  -1.99%  '

# Row 49
$ws.Range("E49").Value = '  +0.63%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.41'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '
